$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dateTimes")

# New row 13: 15th March 2023 "NG-ESO live" saving session (18:30 - 19:00).
$ws.Cells.Item(13, 1).Value = 12
$ws.Cells.Item(13, 2).Value = "NG-ESO live"
$ws.Cells.Item(13, 3).Value = 45000.770833333336
$ws.Cells.Item(13, 4).Value = 45000.791666666664

# Match the date/time number format already used by the column (start_hh/end_hh).
$ws.Range("C13:D13").NumberFormat = "m/d/yy h:mm"
